# Atualização de bases das ligas, do dia: 30-03-2024 às 19:32
#
# 1) Rows 73/74 (A-League fixtures) had their positions swapped in the
#    source feed - everything except the running index in column A moves
#    from row 73 to row 74 and vice versa.
# 2) Rows 133-136 each get a refreshed odds snapshot in place (column A,
#    the running index, never changes) and row 137 (the last data row)
#    is removed outright - net result: the sheet shrinks from 137 to
#    136 data rows (dimension AC137 -> AC136).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Step 1: swap row 73 <-> row 74 (column A keeps its original 71 / 72)
# ---------------------------------------------------------------------

$ws.Range("B73").Value  = 7646749
$ws.Range("C73").Value  = "Australia ALeague"
$ws.Range("D73").Value  = "Australia ALeague"
$ws.Range("E73").Value  = 45305.23958333334
$ws.Range("F73").Value  = "Brisbane Roar"
$ws.Range("G73").Value  = "Newcastle Jets"
$ws.Range("H73").Value  = 3
$ws.Range("I73").Value  = 2
$ws.Range("J73").Value  = "H"
$ws.Range("K73").Value  = 1.909
$ws.Range("L73").Value  = 4
$ws.Range("M73").Value  = 3.4
$ws.Range("N73").Value  = 2.4
$ws.Range("O73").Value  = 4
$ws.Range("P73").Value  = 2.6
$ws.Range("Q73").Value  = 0
$ws.Range("R73").Value  = 1.83
$ws.Range("S73").Value  = 2.07
$ws.Range("T73").Value  = 3.25
$ws.Range("U73").Value  = 1.9
$ws.Range("V73").Value  = 1.95
$ws.Range("W73").Value  = 1.4
$ws.Range("X73").Value  = -1
$ws.Range("Y73").Value  = -1
$ws.Range("Z73").Value  = 0.8300000000000001
$ws.Range("AA73").Value = -1
$ws.Range("AB73").Value = 0.8999999999999999
$ws.Range("AC73").Value = -1

$ws.Range("B74").Value  = 7646750
$ws.Range("C74").Value  = "Australia ALeague"
$ws.Range("D74").Value  = "Australia ALeague"
$ws.Range("E74").Value  = 45305.23958333334
$ws.Range("F74").Value  = "Perth Glory"
$ws.Range("G74").Value  = "Wellington Phoenix"
$ws.Range("H74").Value  = 3
$ws.Range("I74").Value  = 4
$ws.Range("J74").Value  = "A"
$ws.Range("K74").Value  = 2.45
$ws.Range("L74").Value  = 3.75
$ws.Range("M74").Value  = 2.55
$ws.Range("N74").Value  = 3.1
$ws.Range("O74").Value  = 3.8
$ws.Range("P74").Value  = 2.05
$ws.Range("Q74").Value  = 0.25
$ws.Range("R74").Value  = 2
$ws.Range("S74").Value  = 1.85
$ws.Range("T74").Value  = 3
$ws.Range("U74").Value  = 1.925
$ws.Range("V74").Value  = 1.925
$ws.Range("W74").Value  = -1
$ws.Range("X74").Value  = -1
$ws.Range("Y74").Value  = 1.05
$ws.Range("Z74").Value  = -1
$ws.Range("AA74").Value = 0.8500000000000001
$ws.Range("AB74").Value = 0.925
$ws.Range("AC74").Value = -1

# ---------------------------------------------------------------------
# Step 2: refresh the odds feed snapshot for rows 133-135 in place
# (the running index in column A is untouched).
# ---------------------------------------------------------------------

$ws.Range("B133").Value  = 7127394
$ws.Range("E133").Value  = 45381.875
$ws.Range("F133").Value  = "Wellington Phoenix"
$ws.Range("G133").Value  = "Brisbane Roar"
$ws.Range("K133").Value  = 1.8
$ws.Range("L133").Value  = 3.8
$ws.Range("M133").Value  = 4
$ws.Range("N133").Value  = 2.2
$ws.Range("O133").Value  = 3.5
$ws.Range("P133").Value  = 3.2
$ws.Range("Q133").Value  = -0.25
$ws.Range("R133").Value  = 1.98
$ws.Range("S133").Value  = 1.92
$ws.Range("T133").Value  = 2.75
$ws.Range("U133").Value  = 1.825
$ws.Range("V133").Value  = 2.025

$ws.Range("B134").Value  = 7127397
$ws.Range("E134").Value  = 45382.04166666666
$ws.Range("F134").Value  = "Melbourne Victory"
$ws.Range("G134").Value  = "Perth Glory"
$ws.Range("K134").Value  = 1.4
$ws.Range("L134").Value  = 5
$ws.Range("M134").Value  = 6.5
$ws.Range("N134").Value  = 1.4
$ws.Range("O134").Value  = 5.5
$ws.Range("P134").Value  = 6.5
$ws.Range("Q134").Value  = -1.5
$ws.Range("R134").Value  = 2.05
$ws.Range("S134").Value  = 1.85
$ws.Range("T134").Value  = 3.25
$ws.Range("U134").Value  = 1.875
$ws.Range("V134").Value  = 1.975

$ws.Range("B135").Value  = 7127398
$ws.Range("E135").Value  = 45383.04166666666
$ws.Range("F135").Value  = "Macarthur FC"
$ws.Range("G135").Value  = "Western Sydney Wanderers"
$ws.Range("K135").Value  = 2.5
$ws.Range("L135").Value  = 3.5
$ws.Range("M135").Value  = 2.625
$ws.Range("N135").Value  = 2.9
$ws.Range("O135").Value  = 3.8
$ws.Range("P135").Value  = 2.25
$ws.Range("Q135").Value  = 0.25
$ws.Range("R135").Value  = 1.84
$ws.Range("S135").Value  = 2.06
$ws.Range("T135").Value  = 3.25
$ws.Range("U135").Value  = 1.925
$ws.Range("V135").Value  = 1.925

# ---------------------------------------------------------------------
# Step 3: row 136 also gets a refreshed fixture (replacing the old
# 7127397 entry with a brand-new one, id 7898681) - again column A
# (134) is left untouched.
# ---------------------------------------------------------------------

$ws.Range("B136").Value  = 7898681
$ws.Range("E136").Value  = 45384.20833333334
$ws.Range("F136").Value  = "Central Coast Mariners"
$ws.Range("G136").Value  = "Melbourne City"
$ws.Range("K136").Value  = 2.1
$ws.Range("L136").Value  = 4
$ws.Range("M136").Value  = 3
$ws.Range("N136").Value  = 2.05
$ws.Range("O136").Value  = 4
$ws.Range("P136").Value  = 3.1
$ws.Range("Q136").Value  = -0.25
$ws.Range("R136").Value  = 1.83
$ws.Range("S136").Value  = 2.07
$ws.Range("T136").Value  = 3
$ws.Range("U136").Value  = 1.9
$ws.Range("V136").Value  = 1.95

# ---------------------------------------------------------------------
# Step 4: the fixture previously on row 137 (id 7127398, now reused on
# row 135 above) is dropped entirely - this is the only row actually
# removed from the sheet, shrinking the dimension from AC137 to AC136.
# ---------------------------------------------------------------------

$ws.Rows(137).Delete()
